# Weekly price-sheet update: a new daily record was inserted ahead of the
# existing "Femacal de La Calera - Berenjena" series, pushing every
# subsequent record down by one row (the sheet is ordered with the newest
# observations first, oldest ones trailing near the bottom).
#
# Net effect on the worksheet:
#   - Row 139 becomes a brand-new record.
#   - Rows 140-171 each now hold what used to be the row above them
#     (140 <- old139, 141 <- old140, ... 171 <- old170).
#   - Used range grows from A1:R170 to A1:R171.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 139:170 down to 140:171, leaving a blank row 139 ready to
# receive the new record (mirrors Excel's own "Insert Cells -> Shift rows
# down" behaviour, carrying formatting such as the date style along).
$ws.Rows("139:139").Insert()

# Populate the newly opened row 139 with the new observation.
$ws.Range("A139").Value2 = 3
$ws.Range("B139").Value2 = "Femacal de La Calera"
$ws.Range("C139").Value2 = "Coquimbo"
$ws.Range("D139").Value2 = 44504
$ws.Range("E139").Value2 = 5
$ws.Range("F139").Value2 = 100112001
$ws.Range("G139").Value2 = "Berenjena"
$ws.Range("H139").Value2 = "Sin especificar"
$ws.Range("I139").Value2 = "Primera"
$ws.Range("J139").Value2 = 80
$ws.Range("K139").Value2 = 8000
$ws.Range("L139").Value2 = 8500
$ws.Range("M139").Value2 = 8250
$ws.Range("N139").Value2 = "$/caja 60 unidades"
$ws.Range("O139").Value2 = "Región de Arica y Parinacota"
$ws.Range("P139").Value2 = 138
$ws.Range("Q139").Value2 = 60
$ws.Range("R139").Value2 = "Hortaliza"
